$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures ---------------------------------------------
# Valor Mora total (sum of the "Valor Mora" column below) grew because a
# new worker (DAGOBERTO BOLAÑOS ORTIZ) was added to the statement.
$ws.Range("E11").Value = 709730
# Cant. Trabajadores (worker count) went from 4 to 5.
$ws.Range("C13").Value = 5

# --- Make room for the bigger worker table -------------------------------
# Today the detail table occupies rows 16-31 (4 workers x 4 periods) and the
# signature block sits right after it on rows 36-37 (leaving a 4-row gap).
# The new statement has 5 workers x 4 periods = 20 rows (16-35), so insert
# 4 blank rows before the old row 36 to push the signature block down to
# rows 40-41 while keeping everything else intact.
$ws.Rows("32:35").Insert()

# Row 31 used to be the last row of the table and therefore carries a
# special "bottom of table" border style. Once we add more rows it becomes
# just another interior row, so first clone that special border down onto
# the new last row (35) and then restyle row 31 back to a normal interior
# row by copying the formatting of row 30 (an ordinary interior row) onto
# it. Rows 32-34 are brand-new interior rows, so they get row 30's look too.
$ws.Range("B31:J31").Copy($ws.Range("B35:J35"))
$ws.Range("B30:J30").Copy($ws.Range("B31:J31"))
$ws.Range("B30:J30").Copy($ws.Range("B32:J32"))
$ws.Range("B30:J30").Copy($ws.Range("B33:J33"))
$ws.Range("B30:J30").Copy($ws.Range("B34:J34"))

# --- Write the refreshed worker/period detail table ----------------------
# Grouped by worker (4 periods each), ordered: 2003, 2002, 1912, 1911.
$data = @(
    @("CC", "73574401",   "HEMITT ENRIQUE ROCHA CASTRO",   "2003", 35112, 877803),
    @("CC", "73574401",   "HEMITT ENRIQUE ROCHA CASTRO",   "2002", 35112, 877803),
    @("CC", "73574401",   "HEMITT ENRIQUE ROCHA CASTRO",   "1912", 33125, 877803),
    @("CC", "73574401",   "HEMITT ENRIQUE ROCHA CASTRO",   "1911", 33125, 877803),
    @("CC", "73581603",   "REINALDO AVILA RICARDO",        "2003", 35112, 877803),
    @("CC", "73581603",   "REINALDO AVILA RICARDO",        "2002", 35112, 877803),
    @("CC", "73581603",   "REINALDO AVILA RICARDO",        "1912", 46805, 877803),
    @("CC", "73581603",   "REINALDO AVILA RICARDO",        "1911", 46805, 877803),
    @("CC", "1143372467", "YERIS ANDRES NAVARRO MERCADO",  "2003", 35112, 877803),
    @("CC", "1143372467", "YERIS ANDRES NAVARRO MERCADO",  "2002", 35112, 877803),
    @("CC", "1143372467", "YERIS ANDRES NAVARRO MERCADO",  "1912", 33125, 877803),
    @("CC", "1143372467", "YERIS ANDRES NAVARRO MERCADO",  "1911", 33125, 877803),
    @("CC", "1124358444", "DAGOBERTO BOLAÑOS ORTIZ",       "2003", 35112, 877803),
    @("CC", "1124358444", "DAGOBERTO BOLAÑOS ORTIZ",       "2002", 35112, 877803),
    @("CC", "1124358444", "DAGOBERTO BOLAÑOS ORTIZ",       "1912", 33125, 877803),
    @("CC", "1124358444", "DAGOBERTO BOLAÑOS ORTIZ",       "1911", 33125, 877803),
    @("CC", "73574315",   "SANDER ORTEGA MARRUGO",         "2003", 35112, 877803),
    @("CC", "73574315",   "SANDER ORTEGA MARRUGO",         "2002", 35112, 877803),
    @("CC", "73574315",   "SANDER ORTEGA MARRUGO",         "1912", 33125, 877803),
    @("CC", "73574315",   "SANDER ORTEGA MARRUGO",         "1911", 33125, 877803)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row++
}
